# Apply the multiplication-table value updates described by the diff.
$d = $word.ActiveDocument

$replacements = @(
    @("388×3=1164", "455×7=3185"),
    @("881×7=6167", "109×3=327"),
    @("826×3=2478", "859×2=1718"),
    @("506×2=1012", "395×9=3555"),
    @("341×5=1705", "976×2=1952"),
    @("517×6=3102", "599×3=1797"),
    @("784×5=3920", "545×9=4905"),
    @("219×6=1314", "396×9=3564"),
    @("238×2=476",  "920×6=5520"),
    @("856×8=6848", "764×3=2292"),
    @("988×2=1976", "361×2=722"),
    @("772×8=6176", "266×3=798"),
    @("646×7=4522", "347×3=1041"),
    @("486×2=972",  "233×2=466"),
    @("454×8=3632", "213×8=1704"),
    @("339×2=678",  "899×3=2697"),
    @("275×3=825",  "631×3=1893"),
    @("854×8=6832", "811×3=2433"),
    @("580×7=4060", "537×3=1611"),
    @("586×8=4688", "565×8=4520"),
    @("124×8=992",  "827×9=7443"),
    @("879×3=2637", "322×7=2254"),
    @("536×3=1608", "407×7=2849"),
    @("128×5=640",  "426×7=2982"),
    @("433×7=3031", "747×4=2988")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
